# Update cryptos list figures (Price / Volume(1h)) per the Mar 20 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.257.03"
$ws.Range("E2").Value = "  +3.99%  "
$ws.Range("D3").Value = "1.786.97"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'338.56"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.3833"
$ws.Range("E7").Value = "  -3.10%  "
$ws.Range("D8").Value = "'0.3446"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("D9").Value = "'46.84"
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("D11").Value = "'0.07394"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "'23.37"
$ws.Range("E12").Value = "  +8.27%  "
$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "'6.460"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "'7.345"
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("D16").Value = "1.785.62"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "'0.00001077"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "'0.06681"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "'82.26"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "'17.48"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").Value = "'6.447"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "28.265.80"
$ws.Range("E23").Value = "  +4.00%  "
$ws.Range("D24").Value = "'12.08"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Value = "'1.441"
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("D27").Value = "'20.71"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "'2.422"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").Value = "'154.95"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "1.988.46"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "'134.92"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "'4.009"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "'6.122"
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("D34").Value = "'0.08906"
$ws.Range("E34").Value = "  +1.76%  "
$ws.Range("D35").Value = "'12.78"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "'0.02419"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "'0.6868"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").Value = "'5.368"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "'0.06399"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "'0.2167"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "'1.245"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("D42").Value = "'1.497"
$ws.Range("E42").Value = "  -7.03%  "
$ws.Range("D43").Value = "'8.281"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "'14.11"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "'0.9991"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'0.6313"
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "'3.877"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "'133.68"
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("D49").Value = "'2.082"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("D50").Value = "'0.07490"
$ws.Range("E50").Value = "  +5.44%  "
$ws.Range("D51").Value = "'1.209"
$ws.Range("E51").Value = "  +6.35%  "
